$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 201.6439274931327
$ws.Range("H4").Value = 1.586955667810023
$ws.Range("I4").Value = 1109.995995204532
$ws.Range("J4").Value = 0.1541644854914921
$ws.Range("L4").Value = 1.417860270521487
$ws.Range("M4").Value = 0.6599321410685661
